# Update cached market-price / profit figures pulled by the scheduled runner.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1882.909
$ws.Range("I40").Value = 1740.4
$ws.Range("J40").Value = 2001.6666
$ws.Range("K40").Value = 1740.4
$ws.Range("L40").Value = 2001.6666
$ws.Range("M40").Value = -1565.4
$ws.Range("N40").Value = -2351.6666
$ws.Range("H98").Value = 770.5
$ws.Range("I98").Value = 770.5
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 770.5
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 727.5
$ws.Range("N98").ClearContents()
$ws.Range("H122").Value = 770.5
$ws.Range("I122").Value = 770.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2311.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 138.5
$ws.Range("N122").ClearContents()
$ws.Range("H125").Value = 14563.25
$ws.Range("I125").Value = 3733
$ws.Range("J125").Value = 21061.4
$ws.Range("K125").Value = 33597
$ws.Range("L125").Value = 189552.6
$ws.Range("M125").Value = -31137
$ws.Range("N125").Value = -194472.6
$ws.Range("H137").Value = 536470.5600000001
$ws.Range("I137").Value = 1390.6857
$ws.Range("J137").Value = 1638105.6
$ws.Range("K137").Value = 4172.0571
$ws.Range("L137").Value = 4914316.800000001
$ws.Range("M137").Value = -1622.0571
$ws.Range("N137").Value = -4919416.800000001

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5868.533
$ws.Range("I61").Value = 4443.095
$ws.Range("J61").Value = 9194.556
$ws.Range("K61").Value = 4443.095
$ws.Range("L61").Value = 9194.556
$ws.Range("M61").Value = -4231.095
$ws.Range("N61").Value = -9618.556
$ws.Range("H74").Value = 5574.921
$ws.Range("I74").Value = 3137.6333
$ws.Range("J74").Value = 14714.75
$ws.Range("K74").Value = 3137.6333
$ws.Range("L74").Value = 14714.75
$ws.Range("M74").Value = -2263.6333
$ws.Range("N74").Value = -16462.75
$ws.Range("H77").Value = 5574.921
$ws.Range("I77").Value = 3137.6333
$ws.Range("J77").Value = 14714.75
$ws.Range("K77").Value = 15688.1665
$ws.Range("L77").Value = 73573.75
$ws.Range("M77").Value = -11320.1665
$ws.Range("N77").Value = -82309.75
$ws.Range("H122").Value = 4466980
$ws.Range("I122").Value = 3570.75
$ws.Range("J122").Value = 7814537
$ws.Range("K122").Value = 10712.25
$ws.Range("L122").Value = 23443611
$ws.Range("M122").Value = -8262.25
$ws.Range("N122").Value = -23448511
$ws.Range("H136").Value = 5868.533
$ws.Range("I136").Value = 4443.095
$ws.Range("J136").Value = 9194.556
$ws.Range("K136").Value = 13329.285
$ws.Range("L136").Value = 27583.668
$ws.Range("M136").Value = -10779.285
$ws.Range("N136").Value = -32683.668

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 416
$ws.Range("I29").Value = 416
$ws.Range("K29").Value = 416
$ws.Range("M29").Value = -127
$ws.Range("H80").Value = 220.72728
$ws.Range("J80").Value = 236.72223
$ws.Range("L80").Value = 236.72223
$ws.Range("N80").Value = -2232.72223
$ws.Range("H83").Value = 220.72728
$ws.Range("J83").Value = 236.72223
$ws.Range("L83").Value = 1183.61115
$ws.Range("N83").Value = -11167.61115
$ws.Range("H106").Value = 28000
$ws.Range("J106").Value = 28000
$ws.Range("L106").Value = 28000
$ws.Range("N106").Value = -30524
$ws.Range("H118").Value = 57153.875
$ws.Range("J118").Value = 57153.875
$ws.Range("L118").Value = 57153.875
$ws.Range("N118").Value = -60467.875
$ws.Range("H130").Value = 68666.664
$ws.Range("J130").Value = 68666.664
$ws.Range("L130").Value = 68666.664
$ws.Range("N130").Value = -78706.664

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3500
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 3500
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 3500
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -5746
$ws.Range("H89").Value = 3500
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 3500
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 17500
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -28732

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 38937.668
$ws.Range("J64").Value = 38937.668
$ws.Range("L64").Value = 38937.668
$ws.Range("N64").Value = -39433.668
$ws.Range("H67").Value = 38937.668
$ws.Range("J67").Value = 38937.668
$ws.Range("L67").Value = 38937.668
$ws.Range("N67").Value = -40653.668
$ws.Range("H126").Value = 3104.7368
$ws.Range("I126").Value = 1748.75
$ws.Range("K126").Value = 5246.25
$ws.Range("M126").Value = -2776.25

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1016.5
$ws.Range("I22").Value = 1274
$ws.Range("J22").Value = 965
$ws.Range("K22").Value = 1274
$ws.Range("L22").Value = 965
$ws.Range("M22").Value = -979
$ws.Range("N22").Value = -1555
$ws.Range("H27").Value = 1016.5
$ws.Range("I27").Value = 1274
$ws.Range("J27").Value = 965
$ws.Range("K27").Value = 1274
$ws.Range("L27").Value = 965
$ws.Range("M27").Value = -1167
$ws.Range("N27").Value = -1179
$ws.Range("H30").Value = 3901.5
$ws.Range("I30").Value = 1868.6666
$ws.Range("K30").Value = 1868.6666
$ws.Range("M30").Value = -1760.6666
$ws.Range("H35").Value = 21128.572
$ws.Range("I35").Value = 6986.2
$ws.Range("K35").Value = 6986.2
$ws.Range("M35").Value = -6650.2
$ws.Range("H46").Value = 702
$ws.Range("I46").Value = 500
$ws.Range("J46").Value = 742.4
$ws.Range("K46").Value = 500
$ws.Range("L46").Value = 742.4
$ws.Range("M46").Value = -312
$ws.Range("N46").Value = -1118.4
$ws.Range("H55").Value = 445133.44
$ws.Range("I55").Value = 667300.2
$ws.Range("J55").Value = 800
$ws.Range("K55").Value = 667300.2
$ws.Range("L55").Value = 800
$ws.Range("M55").Value = -667127.2
$ws.Range("N55").Value = -1146
$ws.Range("H122").Value = 4703.41
$ws.Range("I122").Value = 4440.6787
$ws.Range("J122").Value = 5372.1816
$ws.Range("K122").Value = 13322.0361
$ws.Range("L122").Value = 16116.5448
$ws.Range("M122").Value = -10872.0361
$ws.Range("N122").Value = -21016.5448

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1552.1428
$ws.Range("I126").Value = 1643.3334
$ws.Range("J126").Value = 1005
$ws.Range("K126").Value = 4930.0002
$ws.Range("L126").Value = 3015
$ws.Range("M126").Value = -2460.0002
$ws.Range("N126").Value = -7955

Write-Output "Applied Pandaemonium Profits update"